# Rename the existing sheet "test1" to "ValidLogin"
$wb = $excel.ActiveWorkbook
$wsValid = $wb.Worksheets.Item(1)
$wsValid.Name = "ValidLogin"

# Add a new worksheet after ValidLogin for the Invalid Login test case
$wsInvalid = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsValid)
$wsInvalid.Name = "InvalidLogin"

# Populate the InvalidLogin sheet with header + invalid credentials
$wsInvalid.Range("A1").Value = "UserName"
$wsInvalid.Range("B1").Value = "Password"
$wsInvalid.Range("A2").Value = "abcd"
$wsInvalid.Range("B2").Value = "xyz"

# Select a cell on the new sheet (mirrors the saved selection state) and
# make it the active tab
$wsInvalid.Range("E7").Select()
$wsInvalid.Activate()
